# "fucked up a merge again....." — untangle the _xlchart defined-name
# indices that got crossed during a bad merge, and restore/append the
# sys2 (sheet7) experiment rows that were dropped.
#
# NOTE: this PowerShell engine does not parse scientific-notation
# numeric literals (e.g. "3.1E-2"), so every literal below is written
# in plain decimal form. It also treats unquoted "$B$1"-style strings
# as variable interpolation, so every sheet reference below is
# single-quoted.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Fix the _xlchart.v1.* defined names whose RefersTo got swapped.
# ---------------------------------------------------------------------
$names = $wb.Names

$names.Item("_xlchart.v1.1").RefersTo  = 'wall_mounted_data!$B$1'
$names.Item("_xlchart.v1.2").RefersTo  = 'wall_mounted_data!$B$2:$B$320'
$names.Item("_xlchart.v1.4").RefersTo  = 'wall_mounted_data!$D$1'
$names.Item("_xlchart.v1.5").RefersTo  = 'wall_mounted_data!$D$2:$D$320'
$names.Item("_xlchart.v1.7").RefersTo  = 'wall_mounted_data!$C$1'
$names.Item("_xlchart.v1.8").RefersTo  = 'wall_mounted_data!$C$2:$C$320'

$names.Item("_xlchart.v1.19").RefersTo = 'adjusted_lens!$D$1'
$names.Item("_xlchart.v1.20").RefersTo = 'adjusted_lens!$D$2:$D$41'
$names.Item("_xlchart.v1.22").RefersTo = 'adjusted_lens!$G$1'
$names.Item("_xlchart.v1.23").RefersTo = 'adjusted_lens!$G$2:$G$41'

# ---------------------------------------------------------------------
# 2. sys2 sheet (sheet7.xml) — re-round a block of values that had
#    picked up spurious long-tail precision, and append the 14 rows
#    of "no-right" classified data that got lost in the merge.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("sys2")

$roundFix = @{
    "F71" = 0.15872112099999999
    "F72" = 0.18191552499999999
    "F73" = 0.221341393
    "F74" = 0.30920677899999999
    "F75" = 0.203403416
    "C76" = 0.21532896700000001
    "F76" = 0.031419546
    "C77" = 0.21251919799999999
    "F77" = 0.027188608
    "C78" = 0.17616757499999999
    "F78" = -0.029476204999999998
    "F79" = 0.27074920400000002
    "F80" = 0.19316558
    "C81" = 0.22792441899999999
    "F81" = 0.076425387999999997
    "C82" = 0.239537202
    "F82" = 0.060881167999999999
    "C83" = 0.171008034
    "F83" = 0.023708579
    "C84" = 0.13651998000000001
    "F84" = -0.00033333400000000002
    "F85" = 0.15508593200000001
    "F86" = 0.20784857700000001
    "C87" = 0.026639979000000001
    "F87" = 0.023169709
}
foreach ($addr in $roundFix.Keys) {
    $ws.Range($addr).Value = $roundFix[$addr]
}

$newRows = @(
    @(0, 0,  0.082672049950766796, -1, -1, 9999,                   1, -1, 1, -1, 0, 0),
    @(0, 0,  0.045515432948655898, -1, -1, 9999,                   1, -1, 1, -1, 0, 0),
    @(0, 0,  0.041364132182983099, -1, -1, 9999,                   1, -1, 1, -1, 0, 0),
    @(0, 0,  0.051255657984232199,  1,  1, 0.054635583596662901,   1,  0, 1,  0, 1, 0),
    @(0, 1, -0.066458714224118795,  1,  1, 0.051641319988196301,   1,  0, 0,  0, 1, 0),
    @(0, 1, -0.105687110478037,     1,  1, -0.016643453998222899,  1,  0, 0,  0, 1, 0),
    @(0, 1, -0.086765498850936096,  1,  1, -0.016101218840755199,  1,  0, 0,  0, 1, 0),
    @(0, 1, -0.077654471799725996,  1,  1, -0.011960002043694301,  1,  0, 0,  0, 1, 0),
    @(0, 1, -0.11111810813152501,   1,  1, -0.011786787797992001,  1,  0, 0,  0, 1, 0),
    @(0, 1, -0.062132519105072302,  1,  0, -0.163323180890215,     1,  0, 0,  1, 1, 0),
    @(0, 1, -0.080074348628915798,  1,  1, -0.0244013273901998,    1,  0, 0,  0, 1, 0),
    @(0, 1, -0.085107871740544405,  1,  1, -0.011141381518335301,  1,  0, 0,  0, 1, 0),
    @(0, 1, -0.086319377561982394,  1,  1, -0.035547655473810801,  1,  0, 0,  0, 1, 0),
    @(0, 1, -0.068487588433220406,  1,  1, -0.021448924428295699,  1,  0, 0,  0, 1, 0)
)

$r = 88
foreach ($row in $newRows) {
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = $row[6]
    $ws.Range("H$r").Value = $row[7]
    $ws.Range("I$r").Value = $row[8]
    $ws.Range("J$r").Value = $row[9]
    $ws.Range("K$r").Value = $row[10]
    $ws.Range("L$r").Value = $row[11]
    $ws.Range("M$r").Value = "no-right"
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 3. Restore the sys2 sheet view: select A2:M101 (whole refreshed
#    range) instead of the stray Q52 selection, and drop the old
#    scrolled-down top-left cell.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2:M101").Select()
